# Weekly driver report update for 2025-04-19
# Updates the "Driver Summary" sheet:
#   - Bad Drivers table (A3:D6) gets refreshed rows (a new worst offender at
#     the top, the rest re-ranked) and new Totals (B7:C7).
#   - Good Drivers table (A15:E.. ) grows from 8 rows to 14 rows: new
#     Intel AX201 driver versions are folded in, the Realtek rows are
#     re-ranked/re-dated, and four brand-new rows are appended at the
#     bottom (23-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Bad Drivers table (header stays put; only the data rows 3-6 and the
#    Totals row 7 change). No structural change, so just overwrite values.
# ---------------------------------------------------------------------

$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.1.2"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 60

$ws.Range("A4").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.155.1"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 32685
$ws.Range("D4").Value = 83.59999999999999

$ws.Range("A5").Value = "Realtek 8821CE Wireless LAN 802.11ac PCI-E NIC - 2024.10.139.3"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 413
$ws.Range("D5").Value = 98.40000000000001

$ws.Range("A6").Value = "Realtek 8821CE Wireless LAN 802.11ac PCI-E NIC - 2024.10.138.0"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 210
$ws.Range("D6").Value = 98.90000000000001

$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 33310

# ---------------------------------------------------------------------
# 2) Good Drivers table grows from rows 15-22 (8 rows) to rows 15-28
#    (14 rows). Clone the existing row formatting down into the six new
#    rows (23-28) first so styles/number formats match, then overwrite
#    every row 15-28 with the final values.
# ---------------------------------------------------------------------

$ws.Range("A22:E22").Copy($ws.Range("A23:E23"))
$ws.Range("A22:E22").Copy($ws.Range("A24:E24"))
$ws.Range("A22:E22").Copy($ws.Range("A25:E25"))
$ws.Range("A22:E22").Copy($ws.Range("A26:E26"))
$ws.Range("A22:E22").Copy($ws.Range("A27:E27"))
$ws.Range("A22:E22").Copy($ws.Range("A28:E28"))

# Column C is an unused spacer in this table; keep the new rows' spacer
# cells as the same blank placeholder text the rest of the table uses.
$ws.Range("C23").Value = "'"
$ws.Range("C24").Value = "'"
$ws.Range("C25").Value = "'"
$ws.Range("C26").Value = "'"
$ws.Range("C27").Value = "'"
$ws.Range("C28").Value = "'"

# Rows 15-18 don't have a known driver-vintage date yet this week, so the
# "Driver Vintage" cell goes back to being an unset numeric cell (0), same
# as the other blank/placeholder cells in column C.
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B15").Value = 56018
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = 0

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B16").Value = 34244
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = 0

$ws.Range("A17").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.143.0"
$ws.Range("B17").Value = 326032
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = 0

$ws.Range("A18").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.144.0"
$ws.Range("B18").Value = 17672
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = 0

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B19").Value = 442178
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").Value = "'2024-11-10"

$ws.Range("A20").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.123.322"
$ws.Range("B20").Value = 16989
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "'2024-06-30"

$ws.Range("A21").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.152.0"
$ws.Range("B21").Value = 1010791
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "'2024-04-15"

$ws.Range("A22").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.149.0"
$ws.Range("B22").Value = 81427
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "'2023-12-20"

$ws.Range("A23").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.141.0"
$ws.Range("B23").Value = 48191
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "'2023-04-17"

$ws.Range("A24").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.128.0"
$ws.Range("B24").Value = 82442
$ws.Range("D24").Value = 99.90000000000001
$ws.Range("E24").Value = "'2022-08-29"

$ws.Range("A25").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.124.0"
$ws.Range("B25").Value = 11789
$ws.Range("D25").Value = 99.90000000000001
$ws.Range("E25").Value = "'2022-07-03"

$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B26").Value = 77849
$ws.Range("D26").Value = 99.90000000000001
$ws.Range("E26").Value = "'2021-08-18"

$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B27").Value = 59673
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = "'2020-08-05"

$ws.Range("A28").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B28").Value = 113652
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = "'2019-12-14"
